# Inventarios - Muestreo Datos.xlsx
# "Creacion de documento Maestro de muestreo de datos, arreglo de muestreos
#  en centro de evento y pagos, y arreglo de diagramas"
#
# This script reworks the sampled "F" (concatenation) formulas on the
# ProductoInventario sheet so they key off Cantidad (C) instead of
# Precio/Fecha (D/E), and then restores each sheet's last-used selection /
# active-tab state the way Excel leaves it after the edits were made.

$wb = $excel.ActiveWorkbook

# --- ProductoInventario: fix the F-column muestreo formula -----------------
# Was: =B2&"-"&D2&"-"&E2   (Producto-Categoria - Precio - Fecha)
# Now: =B2&"-"&C2          (Producto-Categoria - Cantidad)
$wsProductoInventario = $wb.Worksheets.Item("ProductoInventario")
$wsProductoInventario.Activate()

$wsProductoInventario.Range("F2").Formula = "=B2&""-""&C2"
# F3:F4 entered together so Excel records it as one shared formula, same as
# the original authoring (ref="F3:F4" si="0").
$wsProductoInventario.Range("F3:F4").Formula = "=B3&""-""&C3"

$wsProductoInventario.Range("G3").Select() | Out-Null

# --- Inventario: selection left covering the whole used range --------------
$wsInventario = $wb.Worksheets.Item("Inventario")
$wsInventario.Activate()
$wsInventario.Range("A1:D4").Select() | Out-Null

# --- Producto: selection moved back to the header row -----------------------
$wsProducto = $wb.Worksheets.Item("Producto")
$wsProducto.Activate()
$wsProducto.Range("A1:B1").Select() | Out-Null

# --- Sucursal: last sheet touched, stays the active tab on save ------------
$wsSucursal = $wb.Worksheets.Item("Sucursal")
$wsSucursal.Activate()
$wsSucursal.Range("D10").Select() | Out-Null
